$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.733.93"
$ws.Range("D3").Value = "2.477.83"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'320.91"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").Value = "'92.07"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("D11").Value = "'33.02"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "2.859.95"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "'15.49"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").Value = "2.485.56"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "'0.793"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "41.673.24"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "'71.25"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").Value = "'11.24"
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("D23").Value = "'239.58"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("E25").Value = "  +1.45%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'24.97"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'36.59"
$ws.Range("E30").Value = "  +3.62%  "
$ws.Range("D31").Value = "'157.15"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").Value = "'5.43"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'0.0766"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").Value = "'2.57"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").Value = "'17.17"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("E37").Value = "  +2.69%  "
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").Value = "'2.87"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").Value = "'4.01"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D42").Value = "'2.42"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").Value = "1.997.23"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").Value = "'18.69"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("E47").Value = "  +4.61%  "
$ws.Range("D48").Value = "2.738.86"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").Value = "'76.26"
$ws.Range("E49").Value = "  +5.43%  "
$ws.Range("D50").Value = "'97.64"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "'67.26"
$ws.Range("E51").Value = "  +0.30%  "
